$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44161
$ws.Range("J2").Value = 1600
$ws.Range("K2").Value = 1300
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = 1300
$ws.Range("P2").Value = 1300

# Row 3
$ws.Range("D3").Value = 44161
$ws.Range("J3").Value = 1850

# Row 4
$ws.Range("D4").Value = 44162
$ws.Range("J4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44162
$ws.Range("J5").Value = 800

# Row 6
$ws.Range("D6").Value = 44175
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1300
$ws.Range("M6").Value = 1300
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 1300

# Row 7
$ws.Range("D7").Value = 44175
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 1450
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 1000
$ws.Range("P7").Value = 1000

# Row 8
$ws.Range("D8").Value = 44169
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 950
$ws.Range("K8").Value = 1300
$ws.Range("L8").Value = 1300
$ws.Range("M8").Value = 1300
$ws.Range("P8").Value = 1300

# Row 9
$ws.Range("D9").Value = 44169
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 800
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 1000
$ws.Range("P9").Value = 1000

# Row 10
$ws.Range("D10").Value = 44174
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 2800
$ws.Range("K10").Value = 1200
$ws.Range("L10").Value = 1250
$ws.Range("M10").Value = 1221
$ws.Range("P10").Value = 1221

# Row 11
$ws.Range("D11").Value = 44174
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 1300
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 1000
$ws.Range("P11").Value = 1000

# Row 12
$ws.Range("D12").Value = 44172
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 1300
$ws.Range("L12").Value = 1300
$ws.Range("M12").Value = 1300
$ws.Range("P12").Value = 1300

# Row 13
$ws.Range("D13").Value = 44172
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 550
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 1000
$ws.Range("P13").Value = 1000

# Row 14
$ws.Range("D14").Value = 44179
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 980
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 1200
$ws.Range("M14").Value = 1200
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 1200

# Row 15
$ws.Range("D15").Value = 44176
$ws.Range("J15").Value = 2500
$ws.Range("M15").Value = 1256
$ws.Range("P15").Value = 1256

# Row 16
$ws.Range("D16").Value = 44176
$ws.Range("J16").Value = 1500

# Row 17
$ws.Range("D17").Value = 44165
$ws.Range("J17").Value = 720
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 1200
$ws.Range("P17").Value = 1200

# Row 18
$ws.Range("D18").Value = 44165
$ws.Range("J18").Value = 750

# Row 19
$ws.Range("D19").Value = 44159
$ws.Range("J19").Value = 1100

# Row 20
$ws.Range("D20").Value = 44159

# Row 21
$ws.Range("D21").Value = 44168
$ws.Range("J21").Value = 1200

# Row 22
$ws.Range("D22").Value = 44168
$ws.Range("J22").Value = 850

# Row 23
$ws.Range("D23").Value = 44160
$ws.Range("J23").Value = 750
$ws.Range("K23").Value = 1300
$ws.Range("M23").Value = 1300
$ws.Range("P23").Value = 1300

# Row 24
$ws.Range("D24").Value = 44160
$ws.Range("J24").Value = 850

# Row 25
$ws.Range("D25").Value = 44167
$ws.Range("J25").Value = 1430
$ws.Range("K25").Value = 1200
$ws.Range("M25").Value = 1248
$ws.Range("P25").Value = 1248

# Row 26
$ws.Range("D26").Value = 44167
$ws.Range("J26").Value = 350
